$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "y"
$ws.Range("F6").Value = "n"
$ws.Range("F7").Value = "n"
$ws.Range("F8").Value = "n"

$ws.Range("F6:F9").Select()
